$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric need a leading apostrophe
# (quote-prefix) so Excel keeps them as text, exactly like the original
# cells that already contained plain numeric-looking strings.

$ws.Cells.Item(2, 4).Value = '62.771.07'
$ws.Cells.Item(2, 5).Value = '  +6.23%  '

$ws.Cells.Item(3, 4).Value = '3.490.18'
$ws.Cells.Item(3, 5).Value = '  +5.57%  '

$ws.Cells.Item(4, 5).Value = '  +0.39%  '

$ws.Cells.Item(5, 4).Value = "'410.90"
$ws.Cells.Item(5, 5).Value = '  +0.64%  '

$ws.Cells.Item(6, 4).Value = "'129.56"
$ws.Cells.Item(6, 5).Value = '  +15.64%  '

$ws.Cells.Item(7, 4).Value = '3.465.00'
$ws.Cells.Item(7, 5).Value = '  +4.96%  '

$ws.Cells.Item(8, 4).Value = "'0.595"
$ws.Cells.Item(8, 5).Value = '  +2.12%  '

$ws.Cells.Item(9, 5).Value = '  +0.37%  '

$ws.Cells.Item(10, 4).Value = "'0.688"
$ws.Cells.Item(10, 5).Value = '  +9.58%  '

$ws.Cells.Item(11, 4).Value = "'0.128"
$ws.Cells.Item(11, 5).Value = '  +31.05%  '

$ws.Cells.Item(12, 4).Value = "'42.87"
$ws.Cells.Item(12, 5).Value = '  +7.68%  '

$ws.Cells.Item(13, 4).Value = '4.069.77'
$ws.Cells.Item(13, 5).Value = '  +6.45%  '

$ws.Cells.Item(14, 4).Value = "'0.142"
$ws.Cells.Item(14, 5).Value = '  -0.61%  '

$ws.Cells.Item(15, 4).Value = "'8.70"
$ws.Cells.Item(15, 5).Value = '  +2.30%  '

$ws.Cells.Item(16, 4).Value = "'20.01"
$ws.Cells.Item(16, 5).Value = '  +3.56%  '

$ws.Cells.Item(17, 4).Value = '3.454.15'
$ws.Cells.Item(17, 5).Value = '  +3.89%  '

$ws.Cells.Item(18, 4).Value = '62.940.66'
$ws.Cells.Item(18, 5).Value = '  +6.97%  '

$ws.Cells.Item(19, 5).Value = '  +1.01%  '

$ws.Cells.Item(20, 4).Value = "'11.04"
$ws.Cells.Item(20, 5).Value = '  +3.21%  '

$ws.Cells.Item(21, 4).Value = "'0.0000136"
$ws.Cells.Item(21, 5).Value = '  +23.63%  '

$ws.Cells.Item(22, 4).Value = "'3.34"
$ws.Cells.Item(22, 5).Value = '  -0.28%  '

$ws.Cells.Item(23, 4).Value = "'81.65"
$ws.Cells.Item(23, 5).Value = '  +8.81%  '

$ws.Cells.Item(24, 4).Value = "'13.03"
$ws.Cells.Item(24, 5).Value = '  -0.50%  '

$ws.Cells.Item(25, 4).Value = "'311.23"
$ws.Cells.Item(25, 5).Value = '  +1.92%  '

$ws.Cells.Item(26, 4).Value = "'3.16"
$ws.Cells.Item(26, 5).Value = '  -0.36%  '

$ws.Cells.Item(27, 4).Value = "'30.36"
$ws.Cells.Item(27, 5).Value = '  +6.67%  '

$ws.Cells.Item(28, 4).Value = "'8.21"
$ws.Cells.Item(28, 5).Value = '  +5.06%  '

$ws.Cells.Item(29, 2).Value = 'RenderToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(29, 4).Value = "'7.76"
$ws.Cells.Item(29, 5).Value = '  +3.03%  '

$ws.Cells.Item(30, 2).Value = 'Kaspa'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(30, 4).Value = "'0.180"
$ws.Cells.Item(30, 5).Value = '  -0.13%  '

$ws.Cells.Item(31, 4).Value = "'4.37"
$ws.Cells.Item(31, 5).Value = '  -2.20%  '

$ws.Cells.Item(32, 5).Value = '  +2.94%  '

$ws.Cells.Item(33, 2).Value = 'Toncoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(33, 4).Value = "'2.68"
$ws.Cells.Item(33, 5).Value = '  +27.17%  '

$ws.Cells.Item(34, 2).Value = 'Cosmos'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(34, 4).Value = "'11.97"
$ws.Cells.Item(34, 5).Value = '  +4.42%  '

$ws.Cells.Item(35, 4).Value = "'42.90"
$ws.Cells.Item(35, 5).Value = '  +7.83%  '

$ws.Cells.Item(36, 5).Value = '  -0.06%  '

$ws.Cells.Item(37, 4).Value = "'0.0492"
$ws.Cells.Item(37, 5).Value = '  -4.48%  '

$ws.Cells.Item(38, 4).Value = "'52.59"
$ws.Cells.Item(38, 5).Value = '  +1.52%  '

$ws.Cells.Item(39, 5).Value = '  +5.23%  '

$ws.Cells.Item(40, 4).Value = "'0.998"
$ws.Cells.Item(40, 5).Value = '  +0.01%  '

$ws.Cells.Item(41, 4).Value = "'3.01"
$ws.Cells.Item(41, 5).Value = '  -3.75%  '

$ws.Cells.Item(42, 4).Value = "'2.00"
$ws.Cells.Item(42, 5).Value = '  +4.84%  '

$ws.Cells.Item(43, 4).Value = "'137.25"
$ws.Cells.Item(43, 5).Value = '  -1.33%  '

$ws.Cells.Item(44, 4).Value = "'0.125"
$ws.Cells.Item(44, 5).Value = '  +2.40%  '

$ws.Cells.Item(45, 4).Value = "'17.51"
$ws.Cells.Item(45, 5).Value = '  +3.49%  '

$ws.Cells.Item(46, 4).Value = "'0.287"
$ws.Cells.Item(46, 5).Value = '  +1.98%  '

$ws.Cells.Item(47, 4).Value = "'3.96"
$ws.Cells.Item(47, 5).Value = '  +1.06%  '

$ws.Cells.Item(48, 4).Value = "'2.26"
$ws.Cells.Item(48, 5).Value = '  -0.69%  '

$ws.Cells.Item(49, 4).Value = "'22.27"
$ws.Cells.Item(49, 5).Value = '  -0.51%  '

$ws.Cells.Item(50, 4).Value = '2.218.21'
$ws.Cells.Item(50, 5).Value = '  +0.62%  '

$ws.Cells.Item(51, 4).Value = '3.857.75'
$ws.Cells.Item(51, 5).Value = '  +6.04%  '

# Reset style to Normal for the quote-prefixed cells so no stray
# number-format / quotePrefix styling is left behind on those cells.
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(49, 4).Style = "Normal"

